$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: set NumberFormat to Text ("@") before assigning a numeric-looking
# string so Excel keeps it as text, then reset the style to Normal afterwards so no
# stray number-format style index is left on the cell (matches original unstyled cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.451.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.849.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.94%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.850.18"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.84%  "

$ws.Range("E9").Value = "  -1.00%  "

$ws.Range("E10").Value = "  -1.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.47"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.458"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000267"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.498.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.88%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.844.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.553.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.31%  "

$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "471.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.732"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.22%  "

$ws.Range("E24").Value = "  -3.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.14%  "

$ws.Range("E27").Value = "  -1.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.88%  "

$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.002.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.818.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.38%  "

$ws.Range("E38").Value = "  -2.05%  "

$ws.Range("E39").Value = "  -0.88%  "

$ws.Range("E40").Value = "  -2.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.33%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.315"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("E49").Value = "  -2.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.48%  "

# Row 46 becomes Bittensor (was FLOKI)
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "418.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.84%  "

# Row 47 becomes FLOKI (was Bittensor)
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000294"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.72%  "

# Row 51 becomes VeChain (was EnergySwap)
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0359"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.92%  "
